$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): F6 358 -> 359, F7 1776 -> 1778
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F6").Value = 359
$wsExhibit.Range("F7").Value = 1778

# Sheet "全部类型" (sheet4): F6 358 -> 359, F11 1776 -> 1778
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 359
$wsAll.Range("F11").Value = 1778
